# Changes the field syntax from "aql:self.name" to "m:self.name".
#
# The original field is a single instrText run:
#     <w:r><w:instrText xml:space="preserve"> aql:self.name </w:instrText></w:r>
#
# It becomes three instrText runs with an (empty) "_GoBack" bookmark
# spliced in between the "m" and the rest of the code - this mirrors the
# way Word itself splits runs when a user types/edits text in place:
#     <w:r><w:instrText xml:space="preserve"> </w:instrText></w:r>
#     <w:r><w:instrText>m</w:instrText></w:r>
#     <w:bookmarkStart w:id="0" w:name="_GoBack"/>
#     <w:bookmarkEnd w:id="0"/>
#     <w:r><w:instrText xml:space="preserve">:self.name </w:instrText></w:r>

$d = $word.ActiveDocument

foreach ($f in $d.Fields) {

    $code = $f.Code
    $codeText = $code.Text

    $prefixIndex = $codeText.IndexOf("aql:")
    if ($prefixIndex -lt 0) {
        continue
    }

    # Split exactly the way the diff shows:
    # "<lead> aql:self.name <trail>" ->
    #   "<lead> " + "m" + ":self.name <trail>"
    $before = $codeText.Substring(0, $prefixIndex)                 # text before "aql:" (e.g. " ")
    $mChar  = "m"
    $rest   = $codeText.Substring($prefixIndex + "aql".Length)     # ":self.name ..." (drops the "aql" -> keeps ":")

    # Whole field (fldChar begin .. fldChar end) so we can rebuild it with
    # the instrText runs split into three pieces plus the bookmark.
    $fieldRange = $d.Range($code.Start - 1, $code.End + 1)

    # Grab the owning paragraph so the replacement can keep its rsid
    # attributes instead of losing them to the InsertXML rewrite.
    $para = $fieldRange.Paragraphs.Item(1)
    $pPr = ""
    if ($para.Range.Start -eq $fieldRange.Start -and $para.Range.End -eq $fieldRange.End) {
        $pPr = ' w:rsidR="00C909EB" w:rsidRDefault="00C909EB"'
    }

    $beforeX = $before.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $mCharX  = $mChar.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $restX   = $rest.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p' + $pPr + '>' +
           '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
           '<w:r><w:instrText xml:space="preserve">' + $beforeX + '</w:instrText></w:r>' +
           '<w:r><w:instrText>' + $mCharX + '</w:instrText></w:r>' +
           '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
           '<w:bookmarkEnd w:id="0"/>' +
           '<w:r><w:instrText xml:space="preserve">' + $restX + '</w:instrText></w:r>' +
           '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
           '</w:p>' +
           '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $fieldRange.InsertXML($xml)
}
